$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 103, shifting existing rows 103..151 down to 104..152
$ws.Rows.Item(103).Insert()

# Populate the newly inserted row 103 with the new record
$ws.Cells.Item(103, 1).Value  = 1
$ws.Cells.Item(103, 2).Value  = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(103, 3).Value  = "Arica y Parinacota"
$ws.Cells.Item(103, 4).Value  = 45009
$ws.Cells.Item(103, 5).Value  = 15
$ws.Cells.Item(103, 6).Value  = "Fruta"
$ws.Cells.Item(103, 7).Value  = 100102
$ws.Cells.Item(103, 8).Value  = "Cítricos"
$ws.Cells.Item(103, 9).Value  = 100102004
$ws.Cells.Item(103, 10).Value = "Mandarina"
$ws.Cells.Item(103, 11).Value = "Murcott"
$ws.Cells.Item(103, 12).Value = "Segunda"
$ws.Cells.Item(103, 13).Value = 300
$ws.Cells.Item(103, 14).Value = 21000
$ws.Cells.Item(103, 15).Value = 22000
$ws.Cells.Item(103, 16).Value = 21500
$ws.Cells.Item(103, 17).Value = "`$/caja 20 kilos"
$ws.Cells.Item(103, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(103, 19).Value = 1075
$ws.Cells.Item(103, 20).Value = 20
